# Update the "K" column (column G) values on Sheet1 for rows 2-28.
# This corresponds to regenerating the save_data to use K instead of Strike#.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @(1,2,2,2,2,1,3,2,2,1,2,0,2,2,4,1,0,2,2,3,5,1,3,6,5,2,1)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
